$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.103.95"
$ws.Range("E2").Value = "  -1.46%  "

$ws.Range("D3").Value = "1.425.41"
$ws.Range("E3").Value = "  -0.78%  "

$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9969"
$ws.Range("E5").Value = "  -0.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "276.87"
$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3712"
$ws.Range("E7").Value = "  -0.97%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3161"
$ws.Range("E8").Value = "  +3.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.88"
$ws.Range("E9").Value = "  -1.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.064"
$ws.Range("E10").Value = "  +5.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06573"
$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9983"
$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.563"
$ws.Range("E13").Value = "  +4.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.15"
$ws.Range("E14").Value = "  +5.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.232"
$ws.Range("E15").Value = "  +2.07%  "

$ws.Range("D16").Value = "1.426.36"
$ws.Range("E16").Value = "  -1.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001027"
$ws.Range("E17").Value = "  +1.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05730"
$ws.Range("E18").Value = "  -1.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9976"
$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.94"
$ws.Range("E20").Value = "  -5.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.637"
$ws.Range("E21").Value = "  -1.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.90"
$ws.Range("E22").Value = "  +4.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.12"
$ws.Range("E23").Value = "  +3.11%  "

$ws.Range("E24").Value = "  -3.71%  "

$ws.Range("D25").Value = "20.151.42"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.305"
$ws.Range("E26").Value = "  +5.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.90"
$ws.Range("E27").Value = "  -5.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.48"
$ws.Range("E28").Value = "  +3.11%  "

$ws.Range("D29").Value = "1.585.75"
$ws.Range("E29").Value = "  -1.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "111.39"
$ws.Range("E30").Value = "  +1.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.964"
$ws.Range("E31").Value = "  +2.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.332"
$ws.Range("E32").Value = "  -1.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8351"
$ws.Range("E33").Value = "  -7.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07810"
$ws.Range("E34").Value = "  +1.51%  "

$ws.Range("E35").Value = "  +13.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.937"
$ws.Range("E36").Value = "  +5.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05873"
$ws.Range("E37").Value = "  +3.94%  "

$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9966"
$ws.Range("E38").Value = "  -0.59%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.882"
$ws.Range("E39").Value = "  -5.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.78"
$ws.Range("E40").Value = "  -0.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02070"
$ws.Range("E41").Value = "  +1.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.111"
$ws.Range("E42").Value = "  -2.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1878"
$ws.Range("E43").Value = "  -1.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5375"
$ws.Range("E44").Value = "  +1.50%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.44"
$ws.Range("E45").Value = "  +2.64%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.552"
$ws.Range("E46").Value = "  -0.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.09"
$ws.Range("E47").Value = "  +6.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5268"
$ws.Range("E48").Value = "  +2.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.796"
$ws.Range("E49").Value = "  +0.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.044"
$ws.Range("E50").Value = "  -0.63%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06251"
$ws.Range("E51").Value = "  -0.35%  "
